# The deck originally contained a second, unrelated assignment
# ("Class assignment #2" / NYC bike counts) tacked on as slides 10-15.
# This edit removes that trailing assignment, leaving only the original
# 9 "Explanatory analytics" slides (sldId 257-265).
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 10; $i--) {
    $p.Slides.Item($i).Delete()
}

Write-Output "Slide count after delete: $($p.Slides.Count)"
